$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 88
$ws1.Range("F6").Value = 248
$ws1.Range("F7").Value = 199
$ws1.Range("F8").Value = 1909
$ws1.Range("F10").Value = 4463
$ws1.Range("F12").Value = 312

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 49

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 88
$ws4.Range("F5").Value = 49
$ws4.Range("F8").Value = 248
$ws4.Range("F9").Value = 199
$ws4.Range("F12").Value = 1909
$ws4.Range("F14").Value = 4463
$ws4.Range("F16").Value = 312

$wb.Save()
